# Add "investigators" and "funding" columns (Q, R) to the "borehole" sheet,
# with headers, header styling, comments, column widths, and updated
# conditional-formatting formulas that referenced the old last column (P/16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("borehole")

# --- Headers -----------------------------------------------------------
$ws.Range("Q1").Value = "investigators"
$ws.Range("R1").Value = "funding"

# Match the bold / shaded header style used by the rest of row 1 (copy
# from the previously-last header cell, P1).
$ws.Range("Q1:R1").Font.Bold = $true
$ws.Range("Q1:R1").Interior.Color = $ws.Range("P1").Interior.Color

# --- Column widths -------------------------------------------------------
$ws.Columns.Item(17).ColumnWidth = 15.42
$ws.Columns.Item(18).ColumnWidth = 9.75

# --- Comments ------------------------------------------------------------
$investigatorsComment = "[string] investigators`nNames of people and/or agencies who performed the work, as a pipe-delimited list. Each entry should be in the format {person} ({agencies}) [{notes}], where either person or at least one (semicolon-delimited) agencies is required.`nconstraints:`n  - pattern: [^\s]+( [^\s]+)*"
$fundingComment = "[string] funding`nFunding sources as a pipe-delimited list. Each entry should be in the format {funder} [{rorid}] > {award} [{number}] ({url}), where only the funder is required and rorid is the funder's ROR (https://ror.org) ID (e.g. 01jtrvx49).`nconstraints:`n  - pattern: [^\s]+( [^\s]+)*"

$ws.Range("Q1").AddComment($investigatorsComment) | Out-Null
$ws.Range("R1").AddComment($fundingComment) | Out-Null

# --- Conditional formatting ------------------------------------------------
# The formulas for columns A, B, D, E, F reference the blank-count of the
# whole row ($A2:$P2 <> 16); now that the row spans through R, that needs
# to become $A2:$R2 <> 18.
$a = $ws.Range("A2:A1048576").FormatConditions.Item(1)
$a.Formula1 = "=OR(AND(ISBLANK(A2), COUNTBLANK(`$A2:`$R2) <> 18), IF(ISBLANK(A2), FALSE, OR(IF(ISNUMBER(A2), INT(A2) <> A2, TRUE), COUNTIF(A`$2:A`$1048576, A2) >= 2, A2 < 1)))"

$b = $ws.Range("B2:B1048576").FormatConditions.Item(1)
$b.Formula1 = "=AND(ISBLANK(B2), COUNTBLANK(`$A2:`$R2) <> 18)"

$d = $ws.Range("D2:D1048576").FormatConditions.Item(1)
$d.Formula1 = "=OR(AND(ISBLANK(D2), COUNTBLANK(`$A2:`$R2) <> 18), IF(ISBLANK(D2), FALSE, OR(NOT(ISNUMBER(D2)), D2 < -90, D2 > 90)))"

$e = $ws.Range("E2:E1048576").FormatConditions.Item(1)
$e.Formula1 = "=OR(AND(ISBLANK(E2), COUNTBLANK(`$A2:`$R2) <> 18), IF(ISBLANK(E2), FALSE, OR(NOT(ISNUMBER(E2)), E2 < -180, E2 > 180)))"

$f = $ws.Range("F2:F1048576").FormatConditions.Item(1)
$f.Formula1 = "=OR(AND(ISBLANK(F2), COUNTBLANK(`$A2:`$R2) <> 18), IF(ISBLANK(F2), FALSE, OR(NOT(ISNUMBER(F2)), F2 > 9999.0)))"
